$wb = $excel.ActiveWorkbook

# Add the new worksheet "News_Default_Kanal" at the end of the workbook
# (after List_Tab_Menu_Exclusive, the current last sheet).
$lastIndex = $wb.Worksheets.Count
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($lastIndex))
$newSheet.Name = "News_Default_Kanal"

# Header row (bold, like the other "List_*" reference sheets).
$newSheet.Range("A1").Value = "Index"
$newSheet.Range("B1").Value = "Value"
$newSheet.Range("A1:B1").Font.Bold = $true

# Data rows - the leading apostrophe forces Excel to store the index
# column as text (shared string), matching the existing tab sheets.
$newSheet.Range("A2").Value = "'0"
$newSheet.Range("B2").Value = "Berita Utama"

$newSheet.Range("A3").Value = "'1"
$newSheet.Range("B3").Value = "Terkini"

$newSheet.Range("A4").Value = "'2"
$newSheet.Range("B4").Value = "Populer"

# Column B width, sized to fit its contents.
$newSheet.Columns.Item(2).ColumnWidth = 10.85

# Select column C (matches the authored selection state of the sheet).
$newSheet.Columns.Item(3).Select() | Out-Null

# Page setup to match the other worksheets (A4 portrait).
$ps = $newSheet.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Make the new sheet the active / selected tab of the workbook.
$newSheet.Activate() | Out-Null
$excel.ActiveWindow.ScrollWorkbookTabs(16) | Out-Null
